# Bug chase related to nationality
# Nationality was not corrected early enough in the sample
#
# The comparison table's sample counts shift by one row after the fix,
# so three cells need correcting:
#   1. Header acronym "STDT" -> "TSTD"
#   2. "Original goal thwarted (natural hazard)" sample count: 22 -> 21
#   3. "Total" sample count: 374 -> 373
#
# Using Range.Text (rather than Find & Replace) keeps each run's
# xml:space="preserve" attribute intact, matching the original markup.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 1: merged header cell (columns 2-3) - "STDT" -> "TSTD"
$tbl.Cell(1, 2).Range.Text = "TSTD"

# Row 3: "Original goal thwarted (natural hazard)" - Sample column - 22 -> 21
$tbl.Cell(3, 4).Range.Text = "21"

# Row 6: "Total" - Sample column - 374 -> 373
$tbl.Cell(6, 4).Range.Text = "373"
